{"js": "// Template updates for Release 1.0\n//\n// 1) Update the analytical-validity sentence in the \"Test Methodology\"\n//    section: the VAF bins and their average relative-standard-uncertainty\n//    (CV%) values changed.\n// 2) Update the cached \"Reported\" SAVEDATE field result from 25-Oct-2023\n//    to 16-Nov-2023.\n\nconst body = context.document.body;\n\n// --- 1) VAF / CV% sentence -------------------------------------------------\nconst oldVafSentence =\n  \"5%, 10%-20%, 30%-40% and 50% are on average, 10.2%, 10.4%, 3.5% and 4.4\";\nconst newVafSentence =\n  \"2%, 5%-10%, 20%-40% and 50% are on average, 15.4%, 8.6%, 4.0% and 1.8\";\n\nconst vafHits = body.search(oldVafSentence, { matchCase: true });\nvafHits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < vafHits.items.length; i++) {\n  vafHits.items[i].insertText(newVafSentence, \"Replace\");\n}\nawait context.sync();\n\n// --- 2) Reported date --------------------------------------------------\nconst oldDate = \"25-Oct-2023\";\nconst newDate = \"16-Nov-2023\";\n\nconst dateHits = body.search(oldDate, { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < dateHits.items.length; i++) {\n  dateHits.items[i].insertText(newDate, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Template updates for Release 1.0\n#\n# 1) Update the analytical-validity sentence in the \"Test Methodology\"\n#    section: the VAF bins and their average relative-standard-uncertainty\n#    (CV%) values changed.\n# 2) Update the cached \"Reported\" SAVEDATE field result from 25-Oct-2023\n#    to 16-Nov-2023.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n# --- 1) VAF / CV% sentence -------------------------------------------------\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"5%, 10%-20%, 30%-40% and 50% are on average, 10.2%, 10.4%, 3.5% and 4.4\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"2%, 5%-10%, 20%-40% and 50% are on average, 15.4%, 8.6%, 4.0% and 1.8\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find1.Replacement.Text, $wdReplaceAll) | Out-Null\n\n# --- 2) Reported date --------------------------------------------------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"25-Oct-2023\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"16-Nov-2023\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find2.Replacement.Text, $wdReplaceAll) | Out-Null\n"}
